$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.338.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.599.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.26%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "650.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.84%  "
$ws.Range("E7").Value = "  +3.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.403"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.595.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.198"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.285.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.232.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.603.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.484"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "506.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000195"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.780.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.557"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "578.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.39%  "
$ws.Range("E40").Value = "  +5.04%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.918"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +36.54%  "
$ws.Range("E48").Value = "  +4.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0411"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
